$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "acier de <pl>damas</pl>" -> "acier de damas"  (drop the <pl>/</pl>
# tag-marker runs, merge the plain text into a single run)
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("acier de <pl>damas</pl>", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "acier de damas", 2)

# ---------------------------------------------------------------------------
# Change 2: the "<m>cendre</m>" run-pair (end of the "...avecq de la" paragraph)
# becomes "<tl><m>cendre</m></tl>". Each edit is kept inside the boundaries of
# a single existing run so the surrounding run formatting is left untouched.
# ---------------------------------------------------------------------------
$pCendre = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*avecq de la*cendre*") {
        $pCendre = $cand
        break
    }
}

# 2a: the run that contains only "<m>" (immediately before "cendre") becomes "<tl><m>"
$scan = $pCendre.Range.Duplicate
$scan.Find.Execute("<m>cendre") | Out-Null
$openTag = $d.Range($scan.Start, $scan.Start + 3)
$openTag.Text = "<tl><m>"

# 2b: the run that contains only "</m>" (immediately after "cendre") becomes "</m></tl>"
$scan2 = $pCendre.Range.Duplicate
$scan2.Find.Execute("cendre</m>") | Out-Null
$closeTag = $d.Range($scan2.End - 4, $scan2.End)
$closeTag.Text = "</m></tl>"

# ---------------------------------------------------------------------------
# Change 3: the "<m>sable</m>" run-pair (start of the next paragraph) becomes
# "<tl><m>sable</m></tl>".
# ---------------------------------------------------------------------------
$pSable = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "ou *sable*") {
        $pSable = $cand
        break
    }
}

# 3a: the run that contains only "<m>" (immediately before "sable") becomes "<tl><m>"
$scan3 = $pSable.Range.Duplicate
$scan3.Find.Execute("<m>sable") | Out-Null
$openTag2 = $d.Range($scan3.Start, $scan3.Start + 3)
$openTag2.Text = "<tl><m>"

# 3b: the run that contains only "</m>" (immediately after "sable") becomes "</m></tl>"
$scan4 = $pSable.Range.Duplicate
$scan4.Find.Execute("sable</m>") | Out-Null
$closeTag2 = $d.Range($scan4.End - 4, $scan4.End)
$closeTag2.Text = "</m></tl>"

# ---------------------------------------------------------------------------
# Change 4: "<tl><m>pierre de touche</m></tl>" -> "<tl>pierre de touche</tl>"
# (the now-redundant inner "<m>"/"</m>" markers are dropped)
# ---------------------------------------------------------------------------
$pPierre = $pSable
$scan5 = $pPierre.Range.Duplicate
$scan5.Find.Execute("<tl><m>pierre de touche") | Out-Null
$openTag3 = $d.Range($scan5.Start, $scan5.Start + 7)
$openTag3.Text = "<tl>"

$scan6 = $pPierre.Range.Duplicate
$scan6.Find.Execute("pierre de touche</m></tl>") | Out-Null
$closeTag3 = $d.Range($scan6.End - 9, $scan6.End)
$closeTag3.Text = "</tl>"

Write-Host "Done."
